# Weekly update: add a new fruit/vegetable price record as row 534,
# shifting the subsequent rows (534-565) down to (535-566).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 534, pushing everything below it down.
$ws.Rows.Item(534).EntireRow.Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A534").Value = 10
$ws.Range("B534").Value = "Vega Modelo de Temuco"
$ws.Range("C534").Value = "La Araucanía"
$ws.Range("D534").Value = 44585
$ws.Range("E534").Value = 9
$ws.Range("F534").Value = 100114001
$ws.Range("G534").Value = "Papa"
$ws.Range("H534").Value = "Patagonia"
$ws.Range("I534").Value = "1a nueva(o)"
$ws.Range("J534").Value = 480
$ws.Range("K534").Value = 8000
$ws.Range("L534").Value = 8000
$ws.Range("M534").Value = 8000
$ws.Range("N534").Value = "$/saco 25 kilos"
$ws.Range("O534").Value = "Provincia de Cautín"
$ws.Range("P534").Value = 320
$ws.Range("Q534").Value = 25
$ws.Range("R534").Value = "Hortaliza"
